$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most updated cells (coin names/links/URLs, percentage strings, and
# "thousands-dot" prices) are unambiguous text and can be assigned directly.
$plainUpdates = @{
    D2 = "71.428.09"
    E2 = "  +2.33%  "
    D3 = "3.604.89"
    E3 = "  +1.67%  "
    E4 = "  -0.14%  "
    E5 = "  +2.17%  "
    E6 = "  +2.05%  "
    D7 = "3.595.54"
    E7 = "  +1.61%  "
    E8 = "  +0.88%  "
    E9 = "  -0.13%  "
    E10 = "  +17.38%  "
    E11 = "  +2.06%  "
    E12 = "  +1.37%  "
    E13 = "  +7.19%  "
    E14 = "  +0.62%  "
    D15 = "4.173.70"
    E15 = "  +1.46%  "
    D16 = "71.296.46"
    E16 = "  +2.27%  "
    E17 = "  -0.33%  "
    D18 = "3.588.04"
    E18 = "  +0.88%  "
    E19 = "  -0.10%  "
    E20 = "  +6.27%  "
    E21 = "  +0.54%  "
    E22 = "  -1.99%  "
    E23 = "  -11.34%  "
    E24 = "  +2.92%  "
    E25 = "  +5.01%  "
    E26 = "  -0.05%  "
    E27 = "  +1.16%  "
    E28 = "  +0.50%  "
    E29 = "  -0.17%  "
    E30 = "  +2.15%  "
    E31 = "  -2.46%  "
    E32 = "  -0.61%  "
    E33 = "  +0.56%  "
    E34 = "  -0.74%  "
    E35 = "  +5.79%  "
    E36 = "  -3.94%  "
    E37 = "  +4.07%  "
    D38 = "0.0₃0814"
    E38 = "  +5.00%  "
    E39 = "  -1.14%  "
    E40 = "  +0.16%  "
    E41 = "  +3.03%  "
    D42 = "3.538.67"
    E42 = "  +11.15%  "
    E43 = "  +0.90%  "
    E44 = "  +1.49%  "
    E45 = "  +0.48%  "
    B46 = "ThetaToken"
    C46 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    E46 = "  -0.85%  "
    B47 = "ApeXProtocol"
    C47 = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    E47 = "  -2.68%  "
    E48 = "  +1.87%  "
    E49 = "  +2.16%  "
    E50 = "  -0.32%  "
    E51 = "  +0.27%  "
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# These Price cells look like plain decimal numbers, so Excel would silently
# coerce a direct .Value assignment into a floating-point number (and lose
# the exact decimal text, e.g. 586.27 -> 586.26999999999998). Force the cell
# to text format first, assign the literal string, then restore the default
# "Normal" style so the cell ends up styled exactly like its neighbours.
$textUpdates = @{
    D5 = "586.27"
    D6 = "188.50"
    D8 = "0.623"
    D11 = "0.657"
    D12 = "54.80"
    D14 = "9.59"
    D17 = "19.39"
    D19 = "12.43"
    D20 = "568.62"
    D23 = "17.61"
    D24 = "5.07"
    D25 = "4.62"
    D26 = "95.21"
    D27 = "11.44"
    D28 = "2.95"
    D29 = "9.19"
    D30 = "32.51"
    D31 = "7.31"
    D32 = "12.36"
    D33 = "0.116"
    D34 = "64.79"
    D35 = "3.38"
    D36 = "552.64"
    D37 = "0.419"
    D39 = "37.86"
    D41 = "3.30"
    D43 = "3.45"
    D46 = "2.97"
    D47 = "3.47"
    D48 = "9.54"
    D51 = "1.47"
}
foreach ($addr in $textUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$addr]
    $cell.Style = "Normal"
}
